$d = $word.ActiveDocument

# Locate the end of "...geselecteerd." (the text that must remain).
$r1 = $d.Content
$null = $r1.Find.Execute("heeft op die dag. De afspraken hangen af van wat voor klant je hebt geselecteerd.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the end of "...kan worden" (the last bit of text that must be removed).
$r2 = $d.Content
$null = $r2.Find.Execute("op een nette manier zodat het meteen uitgeprint kan worden", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Remove everything between (the two trailing line breaks, the "Print:" heading,
# and the explanatory paragraph that was added by the reverted commit).
$delRange = $d.Range($r1.End, $r2.End)
$delRange.Delete()
